$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.556.07'
$ws.Range('E2').Value = '  +1.73%  '
$ws.Range('D3').Value = '2.037.12'
$ws.Range('E3').Value = '  +2.82%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '257.61'
$ws.Range('E5').Value = '  +5.09%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.625'
$ws.Range('E6').Value = '  -0.58%  '
$ws.Range('B7').Value = 'USDC'
$ws.Range('C7').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('B8').Value = 'Solana'
$ws.Range('C8').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '57.81'
$ws.Range('E8').Value = '  -5.09%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.388'
$ws.Range('E9').Value = '  +1.31%  '
$ws.Range('E10').Value = '  -0.40%  '
$ws.Range('E11').Value = '  -1.47%  '
$ws.Range('E12').Value = '  -0.30%  '
$ws.Range('D13').Value = '2.336.45'
$ws.Range('E13').Value = '  +2.68%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.823'
$ws.Range('E14').Value = '  -2.56%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '21.48'
$ws.Range('E15').Value = '  -2.77%  '
$ws.Range('E16').Value = '  -1.62%  '
$ws.Range('D17').Value = '2.042.45'
$ws.Range('E17').Value = '  +3.32%  '
$ws.Range('D18').Value = '37.494.69'
$ws.Range('E18').Value = '  +1.77%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '70.15'
$ws.Range('E19').Value = '  -0.20%  '
$ws.Range('D20').Value = '0.0₃0857'
$ws.Range('E20').Value = '  -0.47%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.23'
$ws.Range('E21').Value = '  +1.01%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '229.53'
$ws.Range('E22').Value = '  -0.14%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.65'
$ws.Range('E23').Value = '  +5.90%  '
$ws.Range('E24').Value = '  +0.08%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.35'
$ws.Range('E25').Value = '  -1.17%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.15'
$ws.Range('E26').Value = '  -1.32%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '163.73'
$ws.Range('E27').Value = '  +0.26%  '
$ws.Range('E28').Value = '  -5.53%  '
$ws.Range('E29').Value = '  +2.65%  '
$ws.Range('E30').Value = '  +0.33%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.121'
$ws.Range('E31').Value = '  -0.75%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0668'
$ws.Range('E32').Value = '  +7.72%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.75'
$ws.Range('E33').Value = '  -2.26%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.53'
$ws.Range('E34').Value = '  +0.17%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.48'
$ws.Range('E35').Value = '  +9.31%  '
$ws.Range('E36').Value = '  +3.11%  '
$ws.Range('E37').Value = '  -0.18%  '
$ws.Range('E38').Value = '  +2.29%  '
$ws.Range('E39').Value = '  -2.89%  '
$ws.Range('E40').Value = '  +3.89%  '
$ws.Range('E41').Value = '  -2.89%  '
$ws.Range('E42').Value = '  +2.01%  '
$ws.Range('E43').Value = '  +1.09%  '
$ws.Range('B44').Value = 'InjectiveProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.24'
$ws.Range('E44').Value = '  -1.11%  '
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').Value = '1.399.31'
$ws.Range('E45').Value = '  +1.81%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '91.33'
$ws.Range('E46').Value = '  +1.35%  '
$ws.Range('E47').Value = '  +1.46%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.40'
$ws.Range('E48').Value = '  +1.59%  '
$ws.Range('E49').Value = '  +2.11%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.04'
$ws.Range('E50').Value = '  +1.93%  '
$ws.Range('D51').Value = '2.226.77'
$ws.Range('E51').Value = '  +2.70%  '
